$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "status" PASS values that were filled into column F (rows 2-5),
# leaving only the "status" header in F1 (maven/testng run no longer stamps PASS here).
$ws.Range("F2:F5").ClearContents()

# Update the active selection to reflect the now-empty status range.
$ws.Range("F2:F5").Select()
